$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (column C) date for every data row (2-13)
#    from 46066 (2026-02-13) to 46070 (2026-02-17).
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 3).Value = 46070
}

# 2) Swap the full contents of rows 4 and 5 (columns A through Y),
#    since the two cases ("A 53519-2023" and "A 50825-2025") traded places.
$cols4 = @(1,2,4,5,6,7,8,9,10,11,12,13,14,15,16,17)  # A,B,D-Q (plain values)
foreach ($c in $cols4) {
    $tmp = $ws.Cells.Item(4, $c).Value()
    $ws.Cells.Item(4, $c).Value = $ws.Cells.Item(5, $c).Value()
    $ws.Cells.Item(5, $c).Value = $tmp
}
$fcols4 = @(18,19,20,22,23,24,25)  # R (text), S,T,V,W,X,Y (formulas)
foreach ($c in $fcols4) {
    $tmpf = $ws.Cells.Item(4, $c).Formula()
    $ws.Cells.Item(4, $c).Formula = $ws.Cells.Item(5, $c).Formula()
    $ws.Cells.Item(5, $c).Formula = $tmpf
}

# 3) Swap the full contents of rows 10 and 11 (columns A through Q),
#    since the two cases ("A 18968-2025" and "A 53750-2025") traded places.
$cols10 = @(1,2,4,5,6,7,8,9,10,11,12,13,14,15,16,17)  # A,B,D-Q
foreach ($c in $cols10) {
    $tmp = $ws.Cells.Item(10, $c).Value()
    $ws.Cells.Item(10, $c).Value = $ws.Cells.Item(11, $c).Value()
    $ws.Cells.Item(11, $c).Value = $tmp
}
